$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a new value that looks like a plain number need to be
# pre-formatted as Text so Excel keeps storing them as strings (matching
# the source workbook, where every row is an inline string).

$ws.Range("D2").Value = '51.619.93'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.796.45'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.97'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.92'
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.625'
$ws.Range("E9").Value = '  +5.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.82'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0835'
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.93'
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.74'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").Value = '3.228.64'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '2.788.73'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.944'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '51.541.86'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.70'
$ws.Range("E19").Value = '  +3.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.18'
$ws.Range("E20").Value = '  +3.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.51'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '0.0₃0976'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.34'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.76'
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.03'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.165'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.33'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.78'
$ws.Range("E30").Value = '  +5.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.17'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.00'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("E33").Value = '  +9.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0440'
$ws.Range("E34").Value = '  -5.25%  '
$ws.Range("E35").Value = '  -9.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0853'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.88'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.99'
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.12'
$ws.Range("E40").Value = '  -2.70%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("E42").Value = '  -4.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.45'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.69'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").Value = '2.124.34'
$ws.Range("E46").Value = '  +2.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("E47").Value = '  +2.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").Value = '  +5.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.912'
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.38'
$ws.Range("E50").Value = '  -6.31%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  +7.26%  '
